# Update "想去人数" (interest count) values in column F for the data rows
# that changed between the scraper runs. The same update applies to both
# the "展览" sheet and the "全部类型" sheet (they mirror the same data).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1526
    3  = 36
    4  = 986
    5  = 66
    6  = 2423
    8  = 1498
    9  = 70
    10 = 175
    12 = 433
    14 = 15
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
